# Added mtcars data set
#
# Expands the small 3-row (cyl 4/6/8) summary table into the full mtcars
# summary (cyl 4 split by engine V/S into two rows, cyl 6 split into two
# rows, cyl 8 split into two rows but only the V-engine one is populated),
# pushing the footnote row down by one and growing the sheet dimension /
# merged "Cylinder" cells accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert one blank row above the old row 7 so the existing
#    rows 7-11 (and their merges/formatting) slide down to 8-12 with
#    their original per-cell styles intact.
# ---------------------------------------------------------------------
$ws.Rows(7).Insert()

# ---------------------------------------------------------------------
# 2. Format the brand-new row 7 by borrowing formats from cells that
#    already carry the right look after the insert (rows 8-12 kept their
#    original styles, so these sources are never themselves touched).
#    Also fix up E8, whose style needs to change from the "boxed number"
#    look to the plain blank-merge-continuation look.
# ---------------------------------------------------------------------
$ws.Range("E9").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

$ws.Range("F8").Copy() | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null

$ws.Range("G8").Copy() | Out-Null
$ws.Range("G7").PasteSpecial(-4122) | Out-Null

$ws.Range("H8").Copy() | Out-Null
$ws.Range("H7:K7").PasteSpecial(-4122) | Out-Null

$ws.Range("L8").Copy() | Out-Null
$ws.Range("L7").PasteSpecial(-4122) | Out-Null
$ws.Range("L8").Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Write the full mtcars summary values.
# ---------------------------------------------------------------------

# Row 7: cyl=4, engine=0 (V), N=1
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 91
$ws.Range("J7").Value = 2.14

# Row 8: cyl=4, engine=1 (S), N=10 (E8 stays blank - part of the E7:E8 merge)
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10
$ws.Range("H8").Value = 81.8
$ws.Range("I8").Value = 21.87235698318771
$ws.Range("J8").Value = 2.3003
$ws.Range("K8").Value = 0.5982073312080948

# Row 9: cyl=6, engine=0 (V), N=3
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 131.6666666666667
$ws.Range("I9").Value = 37.52776749732568
$ws.Range("J9").Value = 2.755
$ws.Range("K9").Value = 0.1281600561797629

# Row 10: cyl=6, engine=1 (S), N=4 (E10 stays blank - part of the E9:E10 merge)
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4
$ws.Range("H10").Value = 115.25
$ws.Range("I10").Value = 9.178779875342908
$ws.Range("J10").Value = 3.38875
$ws.Range("K10").Value = 0.1162163929916946

# Row 11: cyl=8, engine=0 (V), N=14
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 14
$ws.Range("H11").Value = 209.2142857142857
$ws.Range("I11").Value = 50.97688551827051
$ws.Range("J11").Value = 3.999214285714287
$ws.Range("K11").Value = 0.7594047444769265

# Row 12 (old row 11, the footnote) already has the right values/styles -
# the insert just shifted it down, nothing else to do there.

# ---------------------------------------------------------------------
# 4. Re-create the merge that the row insert doesn't already fix up: the
#    brand-new cyl=4 group needs its own E7:E8 vertical merge (the other
#    merges - E9:E10, E3:K3, E4:K4, H5:I5, J5:K5, E12:K12 - already moved
#    correctly with the inserted row).
# ---------------------------------------------------------------------
$ws.Range("E7:E8").Merge()
